# Add "CharacterDesignId" column to the WeirdLRResults sheet, right before
# the existing "CharacterDesignName" column (old column B), shifting all
# the other header columns one place to the right (old B:U -> new C:V).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WeirdLRResults")

# Insert a new column before column B; this shifts B:U -> C:V.
$ws.Range("B1").EntireColumn.Insert()

# Match the header formatting used by the rest of row 1 (bold font, thin
# box border, centered/top-aligned) before writing the new header text.
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").HorizontalAlignment = -4108
$ws.Range("B1").VerticalAlignment = -4160
$ws.Range("B1").Borders.LineStyle = 1
$ws.Range("B1").Value = "CharacterDesignId"
